$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data (row 2) and set new header row
$ws.Cells.Clear()

$ws.Range("A1").Value = "Time"
$ws.Range("B1").Value = "Player"
$ws.Range("C1").Value = "Coin"

$ws.Range("A1:C1").Font.Bold = $true

$ws.Range("D8").Select()
